$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that look numeric (e.g. "584.66",
# "8.50", "0.0550"). Force those cells to Text format first so Excel
# keeps the exact original text instead of parsing/rounding them as
# floating-point numbers (which would also drop formatting such as
# trailing zeros). Other columns (B/C/E) are never numeric-looking
# (URLs, names, percentages with +/- and spaces) so a plain .Value
# assignment is sufficient and keeps styling untouched.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.222.25'
$ws.Range('E2').Value = '  +0.42%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.574.15'
$ws.Range('E3').Value = '  +1.02%  '

$ws.Range('E4').Value = '  +0.10%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '584.66'
$ws.Range('E5').Value = '  +3.07%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '147.59'
$ws.Range('E6').Value = '  +0.72%  '

$ws.Range('E7').Value = '  +0.04%  '

$ws.Range('E8').Value = '  +3.20%  '

$ws.Range('E9').Value = '  +3.87%  '

$ws.Range('E10').Value = '  +0.67%  '

$ws.Range('E11').Value = '  +0.18%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.357'
$ws.Range('E12').Value = '  +1.63%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '27.53'
$ws.Range('E13').Value = '  +1.27%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.040.09'
$ws.Range('E14').Value = '  +1.34%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '63.236.62'
$ws.Range('E15').Value = '  +0.49%  '

$ws.Range('E16').Value = '  +4.15%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.592.91'
$ws.Range('E17').Value = '  +2.02%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '11.37'
$ws.Range('E18').Value = '  -0.77%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '342.53'
$ws.Range('E19').Value = '  +2.29%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.41'
$ws.Range('E20').Value = '  +3.06%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.87'
$ws.Range('E21').Value = '  +1.79%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.997'
$ws.Range('E22').Value = '  -0.26%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '66.90'
$ws.Range('E23').Value = '  +3.19%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.696.63'
$ws.Range('E24').Value = '  +1.36%  '

$ws.Range('B25').Value = 'Kaspa'
$ws.Range('C25').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.171'
$ws.Range('E25').Value = '  +1.16%  '

$ws.Range('B26').Value = 'Fetch.AI'
$ws.Range('C26').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.64'
$ws.Range('E26').Value = '  +2.17%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.15'
$ws.Range('E27').Value = '  +11.83%  '

$ws.Range('B28').Value = 'Binance-PegBSC-USD'
$ws.Range('C28').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.01'
$ws.Range('E28').Value = '  +1.33%  '

$ws.Range('B29').Value = 'InternetComputer(DFINITY)'
$ws.Range('C29').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.50'
$ws.Range('E29').Value = '  +1.74%  '

$ws.Range('B30').Value = 'SuiNetwork'
$ws.Range('C30').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.49'
$ws.Range('E30').Value = '  -0.59%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.98'
$ws.Range('E31').Value = '  +7.10%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0₃0826'
$ws.Range('E32').Value = '  +2.28%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '462.03'
$ws.Range('E33').Value = '  +13.54%  '

$ws.Range('E34').Value = '  +3.74%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '176.25'
$ws.Range('E35').Value = '  -0.06%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.408'
$ws.Range('E36').Value = '  +2.83%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '19.25'
$ws.Range('E37').Value = '  +1.38%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.53'
$ws.Range('E38').Value = '  +3.81%  '

$ws.Range('E39').Value = '  +0.03%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.74'
$ws.Range('E40').Value = '  -0.76%  '

$ws.Range('E41').Value = '  +0.11%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '151.53'
$ws.Range('E42').Value = '  -0.88%  '

$ws.Range('E43').Value = '  +2.18%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '21.08'
$ws.Range('E44').Value = '  +1.56%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0550'
$ws.Range('E45').Value = '  +6.21%  '

$ws.Range('E46').Value = '  +1.88%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0981'
$ws.Range('E47').Value = '  +2.57%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0239'
$ws.Range('E48').Value = '  +1.73%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.75'
$ws.Range('E49').Value = '  -1.41%  '

$ws.Range('E50').Value = '  -0.08%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.164'
$ws.Range('E51').Value = '  +3.77%  '
